$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 1 de Junio de 2020 a las 18:05'
$ws.Cells.Item(4, 2).Value = 1841698
$ws.Cells.Item(4, 3).Value = 4528
$ws.Cells.Item(4, 4).Value = 600150
$ws.Cells.Item(4, 5).Value = 1135267
$ws.Cells.Item(4, 7).Value = 86
$ws.Cells.Item(4, 8).Value = 106281
$ws.Cells.Item(10, 2).Value = 194837
$ws.Cells.Item(10, 3).Value = 4228
$ws.Cells.Item(10, 5).Value = 95917
$ws.Cells.Item(10, 7).Value = 169
$ws.Cells.Item(10, 8).Value = 5577
$ws.Cells.Item(12, 2).Value = 183596
$ws.Cells.Item(12, 3).Value = 102
$ws.Cells.Item(12, 5).Value = 9086
$ws.Cells.Item(12, 7).Value = 5
$ws.Cells.Item(12, 8).Value = 8610
$ws.Cells.Item(16, 2).Value = 105159
$ws.Cells.Item(16, 3).Value = 5471
$ws.Cells.Item(16, 4).Value = 44946
$ws.Cells.Item(16, 5).Value = 59100
$ws.Cells.Item(16, 7).Value = 59
$ws.Cells.Item(16, 8).Value = 1113
$ws.Cells.Item(17, 2).Value = 91647
$ws.Cells.Item(17, 3).Value = 700
$ws.Cells.Item(17, 4).Value = 49225
$ws.Cells.Item(17, 5).Value = 35097
$ws.Cells.Item(17, 7).Value = 30
$ws.Cells.Item(17, 8).Value = 7325
$ws.Cells.Item(29, 4).Value = 22466
$ws.Cells.Item(29, 5).Value = 12802
$ws.Cells.Item(29, 7).Value = 1
$ws.Cells.Item(29, 8).Value = 24
$ws.Cells.Item(31, 4).Value = 17291
$ws.Cells.Item(31, 5).Value = 16361
$ws.Cells.Item(39, 1).Value = 'Polonia'
$ws.Cells.Item(39, 2).Value = 24165
$ws.Cells.Item(39, 3).Value = 379
$ws.Cells.Item(39, 4).Value = 11449
$ws.Cells.Item(39, 5).Value = 11642
$ws.Cells.Item(39, 8).Value = 1074
$ws.Cells.Item(40, 1).Value = 'Ucrania'
$ws.Cells.Item(40, 2).Value = 24012
$ws.Cells.Item(40, 3).Value = 340
$ws.Cells.Item(40, 4).Value = 9690
$ws.Cells.Item(40, 5).Value = 13604
$ws.Cells.Item(40, 7).Value = 10
$ws.Cells.Item(40, 8).Value = 718
$ws.Cells.Item(43, 2).Value = 17572
$ws.Cells.Item(43, 3).Value = 287
$ws.Cells.Item(43, 4).Value = 10893
$ws.Cells.Item(43, 5).Value = 6177
$ws.Cells.Item(55, 4).Value = 5587
$ws.Cells.Item(55, 5).Value = 5680
$ws.Cells.Item(60, 2).Value = 9286
$ws.Cells.Item(60, 3).Value = 18
$ws.Cells.Item(60, 4).Value = 6642
$ws.Cells.Item(60, 5).Value = 2323
$ws.Cells.Item(60, 7).Value = 1
$ws.Cells.Item(60, 8).Value = 321
$ws.Cells.Item(74, 2).Value = 4019
$ws.Cells.Item(74, 3).Value = 1
$ws.Cells.Item(74, 4).Value = 3845
$ws.Cells.Item(74, 5).Value = 64
$ws.Cells.Item(83, 2).Value = 2918
$ws.Cells.Item(83, 3).Value = 1
$ws.Cells.Item(83, 5).Value = 1365
$ws.Cells.Item(83, 7).Value = 4
$ws.Cells.Item(83, 8).Value = 179
$ws.Cells.Item(92, 2).Value = 2083
$ws.Cells.Item(92, 3).Value = 38
$ws.Cells.Item(92, 4).Value = 1826
$ws.Cells.Item(92, 5).Value = 174
$ws.Cells.Item(93, 1).Value = 'Somalia'
$ws.Cells.Item(93, 2).Value = 2023
$ws.Cells.Item(93, 3).Value = 47
$ws.Cells.Item(93, 4).Value = 361
$ws.Cells.Item(93, 5).Value = 1583
$ws.Cells.Item(93, 7).Value = 1
$ws.Cells.Item(93, 8).Value = 79
$ws.Cells.Item(94, 1).Value = 'Kenia'
$ws.Cells.Item(94, 2).Value = 2021
$ws.Cells.Item(94, 3).Value = 59
$ws.Cells.Item(94, 4).Value = 482
$ws.Cells.Item(94, 5).Value = 1470
$ws.Cells.Item(94, 7).Value = 5
$ws.Cells.Item(94, 8).Value = 69
$ws.Cells.Item(107, 1).Value = 'Mali'
$ws.Cells.Item(107, 2).Value = 1315
$ws.Cells.Item(107, 3).Value = 50
$ws.Cells.Item(107, 4).Value = 744
$ws.Cells.Item(107, 5).Value = 493
$ws.Cells.Item(107, 7).Value = 1
$ws.Cells.Item(107, 8).Value = 78
$ws.Cells.Item(108, 1).Value = 'Guinea Ecuatorial'
$ws.Cells.Item(108, 2).Value = 1306
$ws.Cells.Item(108, 4).Value = 200
$ws.Cells.Item(108, 5).Value = 1094
$ws.Cells.Item(108, 8).Value = 12
$ws.Cells.Item(131, 2).Value = 746
$ws.Cells.Item(131, 3).Value = 7
$ws.Cells.Item(131, 4).Value = 535
$ws.Cells.Item(131, 5).Value = 202
$ws.Cells.Item(153, 2).Value = 296
$ws.Cells.Item(153, 3).Value = 8
$ws.Cells.Item(153, 4).Value = 159
$ws.Cells.Item(153, 5).Value = 110
